$wb = $excel.ActiveWorkbook

# --- Add a new "Unbounded" sheet FIRST (before touching other sheet refs, since
#     inserting a sheet shifts worksheet-collection positions and any sheet
#     object fetched earlier would then resolve against the new layout). ---

$unbounded = $wb.Worksheets.Add()
$unbounded.Name = "Unbounded"

$headers = @("Asset/Liability", "Lower", "Upper")
for ($c = 1; $c -le 3; $c++) {
    $unbounded.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$rows = @(
    @("Liability", -1.0000000000010001, -0.99999999999999001),
    @("15+ STRIPS", 0, 1),
    @("Long Corporate", 0, 1),
    @("Ultra 30Y Futures", 0, 1),
    @("Equity", 0, 1),
    @("Liquid Alternatives", 0, 1),
    @("Private Equity", 0, 1),
    @("Credit", 0, 1),
    @("Real Estate", 0, 1),
    @("Cash", 0.01, 0.0200000000001),
    @("Hedges", 0, 0.5)
)

$r = 2
foreach ($row in $rows) {
    $unbounded.Cells.Item($r, 1).Value = $row[0]
    $unbounded.Cells.Item($r, 2).Value = $row[1]
    $unbounded.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Move it to be the last sheet (after Retirement)
$retirement = $wb.Worksheets.Item("Retirement")
$unbounded.Move($null, $retirement)

# --- Now update bounds on existing sheets (IBT, Pension, Retirement) with fs
#     adjusted weights. Sheet objects are (re-)fetched by name now that the
#     worksheet collection is stable. ---

function Set-Bounds($ws, $data) {
    foreach ($row in $data.Keys) {
        $vals = $data[$row]
        $ws.Cells.Item($row, 2).Value = $vals[0]
        $ws.Cells.Item($row, 3).Value = $vals[1]
    }
}

$ibt = $wb.Worksheets.Item("IBT")
$ibtData = @{
    6  = @(0.15, 0.55000000000000004)
    7  = @(0.03, 0.1)
    8  = @(0.03, 0.1)
    9  = @(0.03, 0.1)
    10 = @(0.03, 0.1)
    11 = @(0.02, 0.0200000000001)
}
Set-Bounds $ibt $ibtData

$pension = $wb.Worksheets.Item("Pension")
$pensionData = @{
    6  = @(0.25, 0.55000000000000004)
    7  = @(0.05, 0.12)
    8  = @(0.05, 0.12)
    9  = @(0.03, 0.1)
    10 = @(0.04, 0.12)
    11 = @(0.02, 0.0200000000001)
}
Set-Bounds $pension $pensionData

$retirement2 = $wb.Worksheets.Item("Retirement")
$retirementData = @{
    6  = @(0.15, 0.35)
    7  = @(0.02, 0.08)
    8  = @(0.02, 0.08)
    9  = @(0.02, 0.08)
    10 = @(0.02, 0.08)
    11 = @(0.02, 0.0200000000001)
}
Set-Bounds $retirement2 $retirementData

# --- Match the final selection state on each sheet ---
$pension.Range("B2:C12").Select() | Out-Null
$retirement2.Range("B2:C12").Select() | Out-Null
$unbounded2 = $wb.Worksheets.Item("Unbounded")
$unbounded2.Range("D9").Select() | Out-Null
$ibt.Range("A9").Select() | Out-Null
